{"js": "// Replace the date and each \"N\u00d7N=\" expression in the document with its\n// updated value, as described by the diff. Every source string in this\n// document is unique, so a simple search+replace per pair is safe and\n// unambiguous.\nconst replacements = [\n  [\"2026-02-04 Wednesday\", \"2026-02-05 Thursday\"],\n  [\"414\u00d74=\", \"364\u00d77=\"],\n  [\"575\u00d79=\", \"171\u00d73=\"],\n  [\"620\u00d72=\", \"748\u00d79=\"],\n  [\"882\u00d79=\", \"716\u00d77=\"],\n  [\"543\u00d74=\", \"318\u00d75=\"],\n  [\"454\u00d72=\", \"169\u00d72=\"],\n  [\"584\u00d72=\", \"675\u00d75=\"],\n  [\"170\u00d77=\", \"201\u00d77=\"],\n  [\"641\u00d75=\", \"568\u00d75=\"],\n  [\"952\u00d73=\", \"746\u00d72=\"],\n  [\"112\u00d78=\", \"520\u00d78=\"],\n  [\"803\u00d75=\", \"588\u00d78=\"],\n  [\"689\u00d72=\", \"699\u00d76=\"],\n  [\"411\u00d75=\", \"128\u00d77=\"],\n  [\"724\u00d74=\", \"792\u00d74=\"],\n  [\"815\u00d77=\", \"878\u00d72=\"],\n  [\"311\u00d78=\", \"860\u00d77=\"],\n  [\"405\u00d79=\", \"161\u00d74=\"],\n  [\"535\u00d78=\", \"880\u00d75=\"],\n  [\"287\u00d73=\", \"328\u00d76=\"],\n  [\"360\u00d77=\", \"289\u00d72=\"],\n  [\"132\u00d72=\", \"520\u00d72=\"],\n  [\"852\u00d76=\", \"484\u00d79=\"],\n  [\"284\u00d75=\", \"503\u00d75=\"],\n  [\"899\u00d79=\", \"983\u00d77=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date and each \"N\u00d7N=\" expression in the document with its\n# updated value, as described by the diff. Every source string in this\n# document is unique, so a simple Find/Replace per pair is safe and\n# unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-02-04 Wednesday\", \"2026-02-05 Thursday\"),\n    @(\"414\u00d74=\", \"364\u00d77=\"),\n    @(\"575\u00d79=\", \"171\u00d73=\"),\n    @(\"620\u00d72=\", \"748\u00d79=\"),\n    @(\"882\u00d79=\", \"716\u00d77=\"),\n    @(\"543\u00d74=\", \"318\u00d75=\"),\n    @(\"454\u00d72=\", \"169\u00d72=\"),\n    @(\"584\u00d72=\", \"675\u00d75=\"),\n    @(\"170\u00d77=\", \"201\u00d77=\"),\n    @(\"641\u00d75=\", \"568\u00d75=\"),\n    @(\"952\u00d73=\", \"746\u00d72=\"),\n    @(\"112\u00d78=\", \"520\u00d78=\"),\n    @(\"803\u00d75=\", \"588\u00d78=\"),\n    @(\"689\u00d72=\", \"699\u00d76=\"),\n    @(\"411\u00d75=\", \"128\u00d77=\"),\n    @(\"724\u00d74=\", \"792\u00d74=\"),\n    @(\"815\u00d77=\", \"878\u00d72=\"),\n    @(\"311\u00d78=\", \"860\u00d77=\"),\n    @(\"405\u00d79=\", \"161\u00d74=\"),\n    @(\"535\u00d78=\", \"880\u00d75=\"),\n    @(\"287\u00d73=\", \"328\u00d76=\"),\n    @(\"360\u00d77=\", \"289\u00d72=\"),\n    @(\"132\u00d72=\", \"520\u00d72=\"),\n    @(\"852\u00d76=\", \"484\u00d79=\"),\n    @(\"284\u00d75=\", \"503\u00d75=\"),\n    @(\"899\u00d79=\", \"983\u00d77=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
